$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "kNN" rating in the CARBON and SILICON
# table (row 20, col B): 94.7 -> 94.8
$ws.Range("B20").Value = 94.8

# Add new row 21: Random Forest Classifier, 95.1
$ws.Range("A21").Value = "Random Forest Classifier"
$ws.Range("B21").Value = 95.1

# Move the selection to reflect where the user ended up after typing
# the new row (matches the post-edit sheetView/selection state)
$ws.Range("B22").Select() | Out-Null
